# Refreshes the cryptocurrency price/volume snapshot on Sheet1 with the latest
# scrape results (GitHub Actions "Updated cryptos list" run). Only the Price
# (D) and Volume(1h) (E) columns change value for most rows; rows 39-40 also
# swap which coin (HuobiToken/Cronos) occupies each rank, so their Coin name
# (B) and Link (C) cells are updated as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking strings (prices such as "248.75" or
# thousand-grouped values such as "2.057.48") that must stay plain text, exactly
# as they were stored before the edit (t="inlineStr"/shared-string, not numeric).
# Forcing the Text number format before assigning the value, then resetting the
# cell style back to "Normal" afterwards, keeps the string type without leaving a
# different visible style/format behind.
$textCells = @{
    'D2' = '37.186.65'
    'D3' = '2.057.48'
    'D5' = '248.75'
    'D6' = '0.666'
    'D7' = '58.83'
    'D9' = '0.387'
    'D10' = '0.0786'
    'D12' = '15.95'
    'D13' = '2.357.22'
    'D14' = '0.839'
    'D15' = '5.78'
    'D16' = '2.050.72'
    'D17' = '18.09'
    'D18' = '37.168.38'
    'D19' = '75.11'
    'D20' = '0.0₃0902'
    'D21' = '5.38'
    'D22' = '238.03'
    'D24' = '2.48'
    'D25' = '2.18'
    'D26' = '169.78'
    'D28' = '20.11'
    'D30' = '4.84'
    'D32' = '0.0621'
    'D34' = '0.0897'
    'D36' = '2.30'
    'D39' = '3.17'
    'D40' = '0.104'
    'D41' = '5.18'
    'D45' = '96.41'
    'D48' = '1.283.58'
    'D49' = '6.87'
    'D50' = '2.241.33'
    'D51' = '3.59'
}

foreach ($cell in $textCells.Keys) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $textCells[$cell]
    $rng.Style = "Normal"
}

# Remaining cells (coin names, links, and the already-text percentage-change
# values) can be assigned directly; none of them are ever misread as numbers.
$plainCells = @{
    'E2' = '  +0.24%  '
    'E3' = '  -0.58%  '
    'E4' = '  -0.08%  '
    'E5' = '  -1.73%  '
    'E6' = '  -1.47%  '
    'E7' = '  -4.18%  '
    'E8' = '  +0.00%  '
    'E9' = '  -0.58%  '
    'E10' = '  -2.17%  '
    'E11' = '  +0.13%  '
    'E12' = '  -2.16%  '
    'E13' = '  -0.59%  '
    'E14' = '  +2.35%  '
    'E15' = '  +5.34%  '
    'E16' = '  -0.91%  '
    'E17' = '  +19.28%  '
    'E18' = '  +0.30%  '
    'E19' = '  +0.67%  '
    'E20' = '  -2.91%  '
    'E21' = '  -1.85%  '
    'E23' = '  +0.01%  '
    'E24' = '  +2.34%  '
    'E25' = '  -6.16%  '
    'E26' = '  -0.26%  '
    'E27' = '  +1.59%  '
    'E28' = '  -1.22%  '
    'E29' = '  -1.02%  '
    'E30' = '  +1.39%  '
    'E31' = '  +3.01%  '
    'E32' = '  -2.48%  '
    'E33' = '  +3.41%  '
    'E34' = '  +0.03%  '
    'E35' = '  -0.07%  '
    'E36' = '  +0.52%  '
    'E37' = '  -0.07%  '
    'E38' = '  -2.34%  '
    'B39' = 'HuobiToken'
    'C39' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'E39' = '  +12.36%  '
    'B40' = 'Cronos'
    'C40' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'E40' = '  -6.50%  '
    'E41' = '  +11.31%  '
    'E42' = '  -2.01%  '
    'E43' = '  -4.13%  '
    'E44' = '  -0.71%  '
    'E45' = '  -2.30%  '
    'E46' = '  -1.93%  '
    'E47' = '  -1.38%  '
    'E48' = '  -1.52%  '
    'E49' = '  -0.66%  '
    'E50' = '  -0.41%  '
    'E51' = '  -18.06%  '
}

foreach ($cell in $plainCells.Keys) {
    $ws.Range($cell).Value = $plainCells[$cell]
}
